$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-06-15 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-06-16 Monday", 2)

$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "341×8=2728"
$tbl.Cell(1,2).Range.Text = "688×6=4128"
$tbl.Cell(1,3).Range.Text = "179×2=358"
$tbl.Cell(1,4).Range.Text = "639×8=5112"
$tbl.Cell(1,5).Range.Text = "833×2=1666"

$tbl.Cell(5,1).Range.Text = "412×6=2472"
$tbl.Cell(5,2).Range.Text = "743×5=3715"
$tbl.Cell(5,3).Range.Text = "983×2=1966"
$tbl.Cell(5,4).Range.Text = "236×2=472"
$tbl.Cell(5,5).Range.Text = "527×6=3162"

$tbl.Cell(10,1).Range.Text = "961×6=5766"
$tbl.Cell(10,2).Range.Text = "657×2=1314"
$tbl.Cell(10,3).Range.Text = "733×7=5131"
$tbl.Cell(10,4).Range.Text = "552×2=1104"
$tbl.Cell(10,5).Range.Text = "884×5=4420"

$tbl.Cell(15,1).Range.Text = "372×8=2976"
$tbl.Cell(15,2).Range.Text = "543×9=4887"
$tbl.Cell(15,3).Range.Text = "607×4=2428"
$tbl.Cell(15,4).Range.Text = "334×6=2004"
$tbl.Cell(15,5).Range.Text = "846×4=3384"

$tbl.Cell(20,1).Range.Text = "896×3=2688"
$tbl.Cell(20,2).Range.Text = "627×2=1254"
$tbl.Cell(20,3).Range.Text = "987×9=8883"
$tbl.Cell(20,4).Range.Text = "914×3=2742"
$tbl.Cell(20,5).Range.Text = "892×8=7136"
